# Adapt the column header formatting to the respective input file names:
# "<field>_old"/"<field>_new" become "<field>_FV2410"/"<field>_FV2504", and
# wrap the data range in an Excel Table ("Table1") with the header row
# frozen, mirroring the regenerated AHB-diff export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the "_old" / "_new" header suffixes to the matching format
#    version names (FV2410 / FV2504).
$headerRenames = @{
    "A1" = "Segmentname_FV2410"
    "B1" = "Segmentgruppe_FV2410"
    "C1" = "Segment_FV2410"
    "D1" = "Datenelement_FV2410"
    "E1" = "Segment ID_FV2410"
    "F1" = "Code_FV2410"
    "G1" = "Qualifier_FV2410"
    "H1" = "Beschreibung_FV2410"
    "I1" = "Bedingungsausdruck_FV2410"
    "J1" = "Bedingung_FV2410"
    "L1" = "Segmentname_FV2504"
    "M1" = "Segmentgruppe_FV2504"
    "N1" = "Segment_FV2504"
    "O1" = "Datenelement_FV2504"
    "P1" = "Segment ID_FV2504"
    "Q1" = "Code_FV2504"
    "R1" = "Qualifier_FV2504"
    "S1" = "Beschreibung_FV2504"
    "T1" = "Bedingungsausdruck_FV2504"
    "U1" = "Bedingung_FV2504"
}

foreach ($addr in $headerRenames.Keys) {
    $ws.Range($addr).Value = $headerRenames[$addr]
}

# 2) Turn the used range into a native Excel Table ("Table1") so the sheet
#    carries an xl/tables/table1.xml part with an AutoFilter, matching the
#    21 renamed columns.
$dataRange = $ws.Range("A1:U77")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"

# 3) Freeze the header row (split below row 1) and leave the selection in
#    the scrollable pane beneath it.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
